$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2.903977394104004, -3.523428678512573, -6.097240924835205),
    @(1.744969725608826, -0.8538760542869568, -4.247409343719482),
    @(4.650478363037109, 11.04139995574951, -3.084672927856445),
    @(-3.62076735496521, -5.12558650970459, 7.287806987762451),
    @(-5.951565742492676, 1.577590107917786, 6.344515323638916),
    @(-1.728857636451721, 8.543354034423828, 3.901864051818848),
    @(1.824864625930786, -0.6525410413742065, -3.505985021591187),
    @(5.772734642028809, -3.56550669670105, -3.791476011276245),
    @(4.579105854034424, 3.013566493988037, -2.666023969650269),
    @(-5.41840124130249, -4.322909832000732, 2.317283153533936),
    @(-11.62250137329102, -2.614225625991821, 11.06896305084228),
    @(7.474494457244873, -2.812897443771362, -3.427688121795654),
    @(1.629921197891235, -8.41618824005127, 3.490139245986938),
    @(8.126436233520508, -4.144477844238281, -0.5264403820037842),
    @(1.081842660903931, -4.525310039520264, -0.875314474105835),
    @(-1.842840909957886, -1.715674996376038, -1.691839694976807),
    @(-15.37276458740234, -16.33429908752441, 3.632352113723755),
    @(7.317900657653809, -9.003682136535645, -1.13044536113739),
    @(0.0980709120631218, -6.593521118164063, 0.4514724016189575),
    @(9.948039054870604, 6.794189929962158, 0.0589224398136138),
    @(-1.0385662317276, -3.36523699760437, 1.002214074134827),
    @(-2.924617052078247, 0.8201870322227478, 1.637111783027649),
    @(-9.886119842529297, -17.29250526428223, -0.4699813723564148),
    @(5.15168571472168, -2.607833862304688, -3.956591844558716),
    @(-0.8223176002502441, -5.908023357391357, 0.7502790689468384),
    @(7.825499057769775, -0.2695784866809845, -1.960819005966186),
    @(2.197174549102783, -1.826995134353638, -2.098770618438721),
    @(-1.679322838783264, 2.35203742980957, 2.066413402557373),
    @(-7.180881023406982, -6.29631233215332, 12.63969612121582),
    @(0.2972753643989563, 7.379552841186523, -0.9200555682182312)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item(2 + $i, 1 + $j).Value = $row[$j]
    }
}